# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Each cell below keeps its original "text" storage type - a leading apostrophe
# forces Excel to store numeric-looking strings (e.g. "1.00", "560.87") as text
# instead of silently coercing them to numbers, matching the source data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "64.160.76"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "3.091.49"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'560.87"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").Value = "'144.71"
$ws.Range("E6").Value = "  +3.25%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.085.94"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").Value = "'6.15"
$ws.Range("E11").Value = "  -5.39%  "
$ws.Range("D12").Value = "'0.471"
$ws.Range("E12").Value = "  +3.82%  "
$ws.Range("D13").Value = "'0.0000228"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").Value = "'35.14"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "3.593.95"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "64.271.31"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").Value = "3.096.36"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").Value = "'6.77"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "'482.45"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'13.98"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").Value = "'0.675"
$ws.Range("D23").Value = "'7.57"
$ws.Range("E23").Value = "  +4.65%  "
$ws.Range("D24").Value = "'13.94"
$ws.Range("E24").Value = "  +10.56%  "
$ws.Range("D25").Value = "'81.23"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("D28").Value = "'8.01"
$ws.Range("E28").Value = "  +1.06%  "
$ws.Range("D29").Value = "'2.07"
$ws.Range("E29").Value = "  +4.03%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("D31").Value = "'26.29"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "'1.14"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "'2.48"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").Value = "'5.60"
$ws.Range("E34").Value = "  -1.91%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'55.84"
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "'6.20"
$ws.Range("E36").Value = "  +3.49%  "
$ws.Range("D37").Value = "'456.18"
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("D38").Value = "'2.99"
$ws.Range("E38").Value = "  +15.94%  "
$ws.Range("D39").Value = "'0.0408"
$ws.Range("E39").Value = "  +2.91%  "
$ws.Range("D40").Value = "'0.0821"
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").Value = "2.983.93"
$ws.Range("E41").Value = "  -2.88%  "
$ws.Range("D42").Value = "'8.26"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  -3.27%  "
$ws.Range("D44").Value = "'27.94"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").Value = "'0.262"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.14"
$ws.Range("E47").Value = "  +3.61%  "
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").Value = "'120.76"
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("D50").Value = "0.0₃0516"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("E51").Value = "  +0.51%  "
